$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(2, 2).Value2 = 0.1423372299900478
$ws.Cells.Item(3, 2).Value2 = 0.1329922600213536
$ws.Cells.Item(4, 2).Value2 = 0.12733173363047
$ws.Cells.Item(5, 2).Value2 = 0.1250444713735561
$ws.Cells.Item(6, 2).Value2 = 0.1246658480548746
$ws.Cells.Item(7, 2).Value2 = 0.127300808057953
$ws.Cells.Item(8, 2).Value2 = 0.1390990346607026
$ws.Cells.Item(9, 2).Value2 = 0.1628501577206265
$ws.Cells.Item(10, 2).Value2 = 0.1806789579369905
$ws.Cells.Item(11, 2).Value2 = 0.1888730601785795
$ws.Cells.Item(12, 2).Value2 = 0.1919880348911107
$ws.Cells.Item(13, 2).Value2 = 0.191316633383849
$ws.Cells.Item(14, 2).Value2 = 0.1891290895252098
$ws.Cells.Item(15, 2).Value2 = 0.1877907248938584
$ws.Cells.Item(16, 2).Value2 = 0.1801451390543036
$ws.Cells.Item(17, 2).Value2 = 0.1754762620766002
$ws.Cells.Item(18, 2).Value2 = 0.1727987240751077
$ws.Cells.Item(19, 2).Value2 = 0.1718935082359394
$ws.Cells.Item(20, 2).Value2 = 0.1759724567382364
$ws.Cells.Item(21, 2).Value2 = 0.1897712971495054
$ws.Cells.Item(22, 2).Value2 = 0.1988598583985066
$ws.Cells.Item(23, 2).Value2 = 0.1940026922636804
$ws.Cells.Item(24, 2).Value2 = 0.1757481063833666
$ws.Cells.Item(25, 2).Value2 = 0.1563586231874297

$ws.Cells.Item(2, 4).Value2 = 0.1647029574918975
$ws.Cells.Item(3, 4).Value2 = 0.1541879627679208
$ws.Cells.Item(4, 4).Value2 = 0.1477354008524827
$ws.Cells.Item(5, 4).Value2 = 0.145106201638697
$ws.Cells.Item(6, 4).Value2 = 0.1446696225261093
$ws.Cells.Item(7, 4).Value2 = 0.1476999423126841
$ws.Cells.Item(8, 4).Value2 = 0.1610762232502196
$ws.Cells.Item(9, 4).Value2 = 0.18736713578609
$ws.Cells.Item(10, 4).Value2 = 0.2067641981097381
$ws.Cells.Item(11, 4).Value2 = 0.2156158898887952
$ws.Cells.Item(12, 4).Value2 = 0.2189725861398699
$ws.Cells.Item(13, 4).Value2 = 0.218249438171398
$ws.Cells.Item(14, 4).Value2 = 0.2158919477450354
$ws.Cells.Item(15, 4).Value2 = 0.214448558057228
$ws.Cells.Item(16, 4).Value2 = 0.2061863434520319
$ws.Cells.Item(17, 4).Value2 = 0.2011254067888046
$ws.Cells.Item(18, 4).Value2 = 0.1982170679023909
$ws.Cells.Item(19, 4).Value2 = 0.1972327746231173
$ws.Cells.Item(20, 4).Value2 = 0.2016638799650821
$ws.Cells.Item(21, 4).Value2 = 0.2165842648432772
$ws.Cells.Item(22, 4).Value2 = 0.2263636561905855
$ws.Cells.Item(23, 4).Value2 = 0.2211413942288516
$ws.Cells.Item(24, 4).Value2 = 0.2014204323703552
$ws.Cells.Item(25, 4).Value2 = 0.1802437554644172

$ws.Cells.Item(2, 5).Value2 = 0.4063152703956945
$ws.Cells.Item(3, 5).Value2 = 0.3540936194290651
$ws.Cells.Item(4, 5).Value2 = 0.322191273672118
$ws.Cells.Item(5, 5).Value2 = 0.3092278998074249
$ws.Cells.Item(6, 5).Value2 = 0.3070774766057127
$ws.Cells.Item(7, 5).Value2 = 0.3220162996601061
$ws.Cells.Item(8, 5).Value2 = 0.388273410560231
$ws.Cells.Item(9, 5).Value2 = 0.5196605242748547
$ws.Cells.Item(10, 5).Value2 = 0.6173291359623647
$ws.Cells.Item(11, 5).Value2 = 0.662064338964413
$ws.Cells.Item(12, 5).Value2 = 0.679052815140011
$ws.Cells.Item(13, 5).Value2 = 0.6753918279009667
$ws.Cells.Item(14, 5).Value2 = 0.6634609994157046
$ws.Cells.Item(15, 5).Value2 = 0.6561594332492291
$ws.Cells.Item(16, 5).Value2 = 0.6144120955137709
$ws.Cells.Item(17, 5).Value2 = 0.5888826873360671
$ws.Cells.Item(18, 5).Value2 = 0.5742273142249417
$ws.Cells.Item(19, 5).Value2 = 0.5692700158555226
$ws.Cells.Item(20, 5).Value2 = 0.5915973591471868
$ws.Cells.Item(21, 5).Value2 = 0.666964029077306
$ws.Cells.Item(22, 5).Value2 = 0.7165036888911374
$ws.Cells.Item(23, 5).Value2 = 0.6900360696645436
$ws.Cells.Item(24, 5).Value2 = 0.5903699883714495
$ws.Cells.Item(25, 5).Value2 = 0.483935313319563

$ws.Cells.Item(2, 6).Value2 = 3.253407090054452
$ws.Cells.Item(3, 6).Value2 = 3.041634659218175
$ws.Cells.Item(4, 6).Value2 = 2.91252819702467
$ws.Cells.Item(5, 6).Value2 = 2.860137794729894
$ws.Cells.Item(6, 6).Value2 = 2.851451484547283
$ws.Cells.Item(7, 6).Value2 = 2.911820758526147
$ws.Cells.Item(8, 6).Value2 = 3.180190043990876
$ws.Cells.Item(9, 6).Value2 = 3.714285372704836
$ws.Cells.Item(10, 6).Value2 = 4.112195600931045
$ws.Cells.Item(11, 6).Value2 = 4.294577246874496
$ws.Cells.Item(12, 6).Value2 = 4.363850481336044
$ws.Cells.Item(13, 6).Value2 = 4.34892174646501
$ws.Cells.Item(14, 6).Value2 = 4.300272132703412
$ws.Cells.Item(15, 6).Value2 = 4.270500458865797
$ws.Cells.Item(16, 6).Value2 = 4.100305143647404
$ws.Cells.Item(17, 6).Value2 = 3.996255405741124
$ws.Cells.Item(18, 6).Value2 = 3.936536980310905
$ws.Cells.Item(19, 6).Value2 = 3.91633900343831
$ws.Cells.Item(20, 6).Value2 = 4.007318318783661
$ws.Cells.Item(21, 6).Value2 = 4.314555925539025
$ws.Cells.Item(22, 6).Value2 = 4.516579067120176
$ws.Cells.Item(23, 6).Value2 = 4.408639266348985
$ws.Cells.Item(24, 6).Value2 = 4.002316461642891
$ws.Cells.Item(25, 6).Value2 = 3.568885083366439

$ws.Cells.Item(2, 7).Value2 = 0.00248067443123344
$ws.Cells.Item(3, 7).Value2 = 0.002487437468431506
$ws.Cells.Item(4, 7).Value2 = 0.00249179418744478
$ws.Cells.Item(5, 7).Value2 = 0.002493621159712442
$ws.Cells.Item(6, 7).Value2 = 0.002493927648558092
$ws.Cells.Item(7, 7).Value2 = 0.002491818617582623
$ws.Cells.Item(8, 7).Value2 = 0.002482964086232053
$ws.Cells.Item(9, 7).Value2 = 0.002467209984874786
$ws.Cells.Item(10, 7).Value2 = 0.002456601874050975
$ws.Cells.Item(11, 7).Value2 = 0.002451982614984631
$ws.Cells.Item(12, 7).Value2 = 0.00245026285501061
$ws.Cells.Item(13, 7).Value2 = 0.002450631929908272
$ws.Cells.Item(14, 7).Value2 = 0.002451840540303783
$ws.Cells.Item(15, 7).Value2 = 0.002452584677905171
$ws.Cells.Item(16, 7).Value2 = 0.002456907887306195
$ws.Cells.Item(17, 7).Value2 = 0.0024596127430065
$ws.Cells.Item(18, 7).Value2 = 0.002461187947932117
$ws.Cells.Item(19, 7).Value2 = 0.002461724631676974
$ws.Cells.Item(20, 7).Value2 = 0.002459322796046395
$ws.Cells.Item(21, 7).Value2 = 0.002451484744679516
$ws.Cells.Item(22, 7).Value2 = 0.002446533693346168
$ws.Cells.Item(23, 7).Value2 = 0.002449160540423408
$ws.Cells.Item(24, 7).Value2 = 0.002459453818123567
$ws.Cells.Item(25, 7).Value2 = 0.002471301097462618

$ws.Cells.Item(2, 11).Value2 = 0.99794825384825
$ws.Cells.Item(3, 11).Value2 = 0.9013982614304155
$ws.Cells.Item(4, 11).Value2 = 0.8431753647372204
$ws.Cells.Item(5, 11).Value2 = 0.8197094231475717
$ws.Cells.Item(6, 11).Value2 = 0.8158284742454498
$ws.Cells.Item(7, 11).Value2 = 0.8428578484565037
$ws.Cells.Item(8, 11).Value2 = 0.9644343140983267
$ws.Cells.Item(9, 11).Value2 = 1.211537328067379
$ws.Cells.Item(10, 11).Value2 = 1.398807645846261
$ws.Cells.Item(11, 11).Value2 = 1.485340500906204
$ws.Cells.Item(12, 11).Value2 = 1.518309000953707
$ws.Cells.Item(13, 11).Value2 = 1.511199613684823
$ws.Cells.Item(14, 11).Value2 = 1.48804878057814
$ws.Cells.Item(15, 11).Value2 = 1.473894543617121
$ws.Cells.Item(16, 11).Value2 = 1.393180158363634
$ws.Cells.Item(17, 11).Value2 = 1.344013521495128
$ws.Cells.Item(18, 11).Value2 = 1.31586009891987
$ws.Cells.Item(19, 11).Value2 = 1.306349242256431
$ws.Cells.Item(20, 11).Value2 = 1.349234308675705
$ws.Cells.Item(21, 11).Value2 = 1.494843243707919
$ws.Cells.Item(22, 11).Value2 = 1.591178645985963
$ws.Cells.Item(23, 11).Value2 = 1.539652961122101
$ws.Cells.Item(24, 11).Value2 = 1.34687363896137
$ws.Cells.Item(25, 11).Value2 = 1.143715630494683

$wb.Save()
